# Homepage / NewCustomer page
# - rename "Sheet2" -> "NewCustomer"
# - populate the New-Customer test-data table (header + 3 rows)
# - add mailto: hyperlinks on the two e-mail cells that carry one,
#   and mirror the Hyperlink cell style on the third e-mail cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "NewCustomer"

# Column widths (approximate best-fit widths from the original workbook)
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 33
$ws.Columns.Item(4).ColumnWidth = 18.833333333333332
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 28.5
$ws.Columns.Item(9).ColumnWidth = 8.666666666666666

# Header row
$ws.Range("A1").Value = "CustomerName"
$ws.Range("B1").Value = "Gender"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "City"
$ws.Range("E1").Value = "State"
$ws.Range("F1").Value = "Pin"
$ws.Range("G1").Value = "Mobile Number"
$ws.Range("H1").Value = "E-mail"
$ws.Range("I1").Value = "Password"

# Row 2 - Harish
$ws.Range("A2").Value = "Harish"
$ws.Range("D2").Value = "Thiruvananthapuram"
$ws.Range("E2").Value = "Kerala"
$ws.Range("F2").Value = 695002
$ws.Range("G2").Value = 9851478957
$ws.Range("H2").Value = "abcd@gmail.com"
$ws.Range("I2").Value = "Test#123"

# Row 3 - Pavithra
$ws.Range("A3").Value = "Pavithra"
$ws.Range("D3").Value = "Hyderabad"
$ws.Range("E3").Value = "Telungana"
$ws.Range("F3").Value = 500004
$ws.Range("G3").Value = 9851478000
$ws.Range("H3").Value = "efgh@gmail.com"
$ws.Range("I3").Value = "Test#456"

# Row 4 - Dormula
$ws.Range("A4").Value = "Dormula"
$ws.Range("D4").Value = "Madurai"
$ws.Range("E4").Value = "Tamil Nadu"
$ws.Range("F4").Value = 590684
$ws.Range("G4").Value = 9851478123
$ws.Range("H4").Value = "krishnanmeena363@gmail.com"
$ws.Range("I4").Value = "Test#789"

# Gender column
$ws.Range("B2").Value = "Male"
$ws.Range("B3").Value = "Female"
$ws.Range("B4").Value = "Female"

# Address column
$ws.Range("C2").Value = "Swami Vivekananda Lane Karamana"
$ws.Range("C3").Value = "Narayana Kandalu Vedhiranth"
$ws.Range("C4").Value = "Alagar nagar Pudhur"

# E-mail hyperlinks
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:abcd@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:efgh@gmail.com")
$ws.Range("H4").Style = "Hyperlink"
